# FAIR4FOIS Research Resource Metadata Schema - vocabulary update
#
# Commit message: "Updating vocabulary and adjusting schema templates"
#
# Content changes (Research_Resource sheet, example row 3):
#   - Contributors (C3): the two ROR identifiers were separated by a plain
#     space; a comma was added between them.
#   - Keywords (H3): the keyword list had a trailing ", " which was removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Research_Resource")

$ws.Range("C3").Value = "https://ror.org/05wg1m734, https://ror.org/02vjkv261"
$ws.Range("H3").Value = "rare diseases, FAIR data, resource discovery"

# The author's selection/cursor ended up on H4 after editing the keywords
# cell next to it.
$ws.Range("H4").Select()
